# dark font for vertex labels
#
# The graph diagram's vertex "type" labels (Account / Person / Company),
# drawn as Oval shapes styled with the accent1 fillRef, currently render
# their text in schemeClr "bg1" (white). Change those labels to the dark
# srgbClr 312D2A used elsewhere in the deck (e.g. edge-label boxes),
# matching the commit "dark font for vertex labels".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Dark color 312D2A as a COM RGB long: R + G*256 + B*65536
$darkColor = 0x31 + (0x2D * 256) + (0x2A * 65536)

# Shape indices (1-based) of the vertex-type label Ovals on slide 1 that
# use <a:schemeClr val="bg1"/> for their text today:
#   4  -> Oval 66  (Account)
#   6  -> Oval 86  (Person)
#   7  -> Oval 91  (Person)
#   21 -> Oval 56  (Account)
#   28 -> Oval 33  (Person)
#   29 -> Oval 115 (Company)
#   37 -> Oval 32  (Account)
#   40 -> Oval 31  (Account)
$vertexLabelIndices = @(4, 6, 7, 21, 28, 29, 37, 40)

foreach ($idx in $vertexLabelIndices) {
    $shp = $s.Shapes.Item($idx)
    $shp.TextFrame.TextRange.Font.Color.RGB = $darkColor
}
